# Update workbook: "Förändrad" (column C) date moves from 45202 -> 45203
# for every existing data row, and two new rows (348, 349) are appended
# with fresh "A 47297-2023" / "A 47326-2023" entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastExistingRow = 347

# 1) Bump column C ("Förändrad") from 45202 to 45203 for every data row.
for ($r = 2; $r -le $lastExistingRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45203
}

# 2) Row 347 gains an explicit row height (ht="15" customHeight="1").
$ws.Rows.Item(347).RowHeight = 15

# 3) Append the two new rows at the bottom of the sheet.
#    (Use positional parameters here - named parameters aren't handled
#    reliably by this PowerShell host.)
function Set-DataRow($RowNumber, $Beteckning, $Datum, $Forandrad, $Lan, $Kommun, $Markagare, $Area) {
    $ws.Cells.Item($RowNumber, 1).Value = $Beteckning
    $ws.Cells.Item($RowNumber, 2).Value = $Datum
    $ws.Cells.Item($RowNumber, 2).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($RowNumber, 3).Value = $Forandrad
    $ws.Cells.Item($RowNumber, 3).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($RowNumber, 4).Value = $Lan
    $ws.Cells.Item($RowNumber, 5).Value = $Kommun
    $ws.Cells.Item($RowNumber, 6).Value = $Markagare
    $ws.Cells.Item($RowNumber, 7).Value = $Area
    $ws.Cells.Item($RowNumber, 8).Value = 0
    $ws.Cells.Item($RowNumber, 9).Value = 0
    $ws.Cells.Item($RowNumber, 10).Value = 0
    $ws.Cells.Item($RowNumber, 11).Value = 0
    $ws.Cells.Item($RowNumber, 12).Value = 0
    $ws.Cells.Item($RowNumber, 13).Value = 0
    $ws.Cells.Item($RowNumber, 14).Value = 0
    $ws.Cells.Item($RowNumber, 15).Value = 0
    $ws.Cells.Item($RowNumber, 16).Value = 0
    $ws.Cells.Item($RowNumber, 17).Value = 0
    $ws.Cells.Item($RowNumber, 18).Value = ""
    $ws.Cells.Item($RowNumber, 18).WrapText = $true
}

Set-DataRow 348 "A 47297-2023" 45202 45203 "VÄSTERBOTTENS LÄN" "NORSJÖ" "Holmen skog AB" 4.5
$ws.Rows.Item(348).RowHeight = 15

Set-DataRow 349 "A 47326-2023" 45202 45203 "VÄSTERBOTTENS LÄN" "NORSJÖ" "Holmen skog AB" 2.3
